$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.130.19'
$ws.Range('E2').Value = '  -0.65%  '
$ws.Range('D3').Value = '1.857.28'
$ws.Range('E3').Value = '  -0.67%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '232.75'
$ws.Range('E5').Value = '  -0.97%  '
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4684'
$ws.Range('E7').Value = '  -0.51%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '42.66'
$ws.Range('E8').Value = '  -0.59%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2828'
$ws.Range('E9').Value = '  -1.43%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06444'
$ws.Range('E10').Value = '  -1.97%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.74'
$ws.Range('E11').Value = '  -4.06%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07666'
$ws.Range('E12').Value = '  -4.39%  '
$ws.Range('D13').Value = '1.856.19'
$ws.Range('E13').Value = '  -0.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '93.28'
$ws.Range('E14').Value = '  -3.82%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.034'
$ws.Range('E15').Value = '  -1.55%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.6755'
$ws.Range('E16').Value = '  -1.26%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '266.81'
$ws.Range('E17').Value = '  -0.97%  '
$ws.Range('D18').Value = '30.100.87'
$ws.Range('E18').Value = '  -0.72%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.27'
$ws.Range('E19').Value = '  -5.45%  '
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000007494'
$ws.Range('E21').Value = '  -1.73%  '
$ws.Range('D22').Value = '2.091.22'
$ws.Range('E22').Value = '  -0.87%  '
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('E24').Value = '  -3.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '6.060'
$ws.Range('E25').Value = '  -2.52%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.237'
$ws.Range('E26').Value = '  -2.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.11'
$ws.Range('E27').Value = '  -2.21%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.48'
$ws.Range('E28').Value = '  -2.29%  '
$ws.Range('E29').Value = '  -3.91%  '
$ws.Range('E30').Value = '  -0.26%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09779'
$ws.Range('E31').Value = '  -1.59%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.445'
$ws.Range('E32').Value = '  -1.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.174'
$ws.Range('E33').Value = '  -4.42%  '
$ws.Range('E34').Value = '  -2.83%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.04625'
$ws.Range('E35').Value = '  -1.77%  '
$ws.Range('E36').Value = '  -2.89%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6805'
$ws.Range('E37').Value = '  -2.90%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.712'
$ws.Range('E38').Value = '  +0.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01810'
$ws.Range('E39').Value = '  -3.41%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.707'
$ws.Range('E40').Value = '  +2.69%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.262'
$ws.Range('E41').Value = '  -0.62%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '70.02'
$ws.Range('E42').Value = '  -2.57%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.000'
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8268'
$ws.Range('E44').Value = '  -1.75%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '102.01'
$ws.Range('E45').Value = '  -0.80%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.864'
$ws.Range('E46').Value = '  -5.24%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4025'
$ws.Range('E47').Value = '  -3.42%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.105'
$ws.Range('E48').Value = '  -0.87%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '920.60'
$ws.Range('E49').Value = '  +0.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.863'
$ws.Range('E51').Value = '  -1.35%  '
